$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.333.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.76%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.802.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.63%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.08"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4455"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +10.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3723"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +7.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.76"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.147"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07510"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.53"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.000"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.688"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.288"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.799.31"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001093"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06790"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.73"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.82%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.322"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.322.89"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.81"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.413"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.23%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.02"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.356"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.003.62"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.32"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.243"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.010"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.818"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09333"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2298"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.61%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06327"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02323"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6573"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.162"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.63%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.212"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.53%  "

$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.458"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.177"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.97"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6066"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.791"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.29"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.034"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.156"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.90%  "
